$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 315122
$ws.Range("D2").Value = 401744408
$ws.Range("C4").Value = 314
$ws.Range("D4").Value = 449207
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 14676
$ws.Range("C8").Value = 847
$ws.Range("D8").Value = 1245908
$ws.Range("C10").Value = 115809
$ws.Range("D10").Value = 169704179
$ws.Range("C12").Value = 58409
$ws.Range("D12").Value = 84305345
$ws.Range("C16").Value = 3949
$ws.Range("D16").Value = 5605861
$ws.Range("C20").Value = 6433
$ws.Range("D20").Value = 8977335
$ws.Range("C22").Value = 76223
$ws.Range("D22").Value = 95141217
$ws.Range("C28").Value = 32149
$ws.Range("D28").Value = 47069041
$ws.Range("C30").Value = 11329
$ws.Range("D30").Value = 16295486
$ws.Range("C33").Value = 1547
$ws.Range("D33").Value = 2171807
$ws.Range("C35").Value = 1766
$ws.Range("D35").Value = 2492023
$ws.Range("C36").Value = 95788
$ws.Range("D36").Value = 120673422
$ws.Range("C44").Value = 44000
$ws.Range("D44").Value = 64490087
$ws.Range("C46").Value = 9003
$ws.Range("D46").Value = 12923224
$ws.Range("C48").Value = 1386
$ws.Range("D48").Value = 1925100
$ws.Range("C51").Value = 2241
$ws.Range("D51").Value = 3125564
$ws.Range("C52").Value = 67967
$ws.Range("D52").Value = 85310896
$ws.Range("C58").Value = 27852
$ws.Range("D58").Value = 40850210
$ws.Range("C61").Value = 10923
$ws.Range("D61").Value = 15793057
$ws.Range("C63").Value = 1343
$ws.Range("D63").Value = 1876789
$ws.Range("C67").Value = 1426
$ws.Range("D67").Value = 1995599
$ws.Range("C69").Value = 20173
$ws.Range("D69").Value = 26423774
$ws.Range("C73").Value = 7486
$ws.Range("D73").Value = 10959496
$ws.Range("C75").Value = 5029
$ws.Range("D75").Value = 7301706
$ws.Range("C76").Value = 483
$ws.Range("D76").Value = 683239
$ws.Range("C78").Value = 138462
$ws.Range("D78").Value = 172733006
$ws.Range("C79").Value = 67
$ws.Range("D79").Value = 80285
$ws.Range("C84").Value = 62874
$ws.Range("D84").Value = 92157458
$ws.Range("C86").Value = 10
$ws.Range("D86").Value = 15000
$ws.Range("C87").Value = 29270
$ws.Range("D87").Value = 42346144
$ws.Range("C89").Value = 2699
$ws.Range("D89").Value = 3886643
$ws.Range("C90").Value = 2734
$ws.Range("D90").Value = 3862665
$ws.Range("C91").Value = 31773
$ws.Range("D91").Value = 43055768
$ws.Range("C95").Value = 7746
$ws.Range("D95").Value = 11388970
$ws.Range("C97").Value = 7039
$ws.Range("D97").Value = 10204733
$ws.Range("C99").Value = 514
$ws.Range("D99").Value = 730405
$ws.Range("C100").Value = 480
$ws.Range("D100").Value = 692443
$ws.Range("C101").Value = 8618
$ws.Range("D101").Value = 11954383
$ws.Range("C103").Value = 2177
$ws.Range("D103").Value = 3207470
$ws.Range("C105").Value = 2921
$ws.Range("D105").Value = 4265402
$ws.Range("C107").Value = 122
$ws.Range("D107").Value = 177120
$ws.Range("C108").Value = 163
$ws.Range("D108").Value = 230586
$ws.Range("C109").Value = 138968
$ws.Range("D109").Value = 171882772
$ws.Range("C115").Value = 52152
$ws.Range("D115").Value = 76458543
$ws.Range("C117").Value = 26529
$ws.Range("D117").Value = 38432803
$ws.Range("C118").Value = 1297
$ws.Range("D118").Value = 1775051
$ws.Range("C121").Value = 2184
$ws.Range("D121").Value = 3067697
$ws.Range("C123").Value = 491420
$ws.Range("D123").Value = 648119682
$ws.Range("C128").Value = 1357
$ws.Range("D128").Value = 2011811
$ws.Range("C129").Value = 31
$ws.Range("D129").Value = 40510
$ws.Range("C130").Value = 204097
$ws.Range("D130").Value = 300042477
$ws.Range("C133").Value = 176257
$ws.Range("D133").Value = 256200655
$ws.Range("C136").Value = 2799
$ws.Range("D136").Value = 3933784
$ws.Range("C138").Value = 6131
$ws.Range("D138").Value = 8661949
$ws.Range("C141").Value = 43635
$ws.Range("D141").Value = 58269634
$ws.Range("C146").Value = 452
$ws.Range("D146").Value = 677665
$ws.Range("C147").Value = 13859
$ws.Range("D147").Value = 20328132
$ws.Range("C148").Value = 3688
$ws.Range("D148").Value = 5319238
$ws.Range("C151").Value = 392
$ws.Range("D151").Value = 563931
$ws.Range("C153").Value = 372
$ws.Range("D153").Value = 524751
$ws.Range("C154").Value = 17170
$ws.Range("D154").Value = 22695489
$ws.Range("C158").Value = 7026
$ws.Range("D158").Value = 10219841
$ws.Range("C160").Value = 4894
$ws.Range("D160").Value = 7044436
$ws.Range("C162").Value = 269
$ws.Range("D162").Value = 371235
$ws.Range("C163").Value = 259
$ws.Range("D163").Value = 370774
$ws.Range("C165").Value = 15080
$ws.Range("D165").Value = 21883231
$ws.Range("C166").Value = 1724
$ws.Range("D166").Value = 2564530
$ws.Range("C167").Value = 233
$ws.Range("D167").Value = 344302
$ws.Range("C170").Value = 78
$ws.Range("D170").Value = 116949
$ws.Range("C171").Value = 86386
$ws.Range("D171").Value = 108082934
$ws.Range("C178").Value = 33500
$ws.Range("D178").Value = 49130831
$ws.Range("C180").Value = 12807
$ws.Range("D180").Value = 18504388
$ws.Range("C182").Value = 1236
$ws.Range("D182").Value = 1729896
$ws.Range("C184").Value = 1598
$ws.Range("D184").Value = 2249693
$ws.Range("C186").Value = 234507
$ws.Range("D186").Value = 291585336
$ws.Range("C188").Value = 165
$ws.Range("D188").Value = 237736
$ws.Range("C194").Value = 85725
$ws.Range("D194").Value = 125673111
$ws.Range("C197").Value = 32542
$ws.Range("D197").Value = 46835537
$ws.Range("C200").Value = 5027
$ws.Range("D200").Value = 7164262
$ws.Range("C203").Value = 4706
$ws.Range("D203").Value = 6509827
$ws.Range("C206").Value = 259210
$ws.Range("D206").Value = 320863622
$ws.Range("C207").Value = 156
$ws.Range("D207").Value = 170973
$ws.Range("C208").Value = 249
$ws.Range("D208").Value = 356087
$ws.Range("C213").Value = 607
$ws.Range("D213").Value = 884406
$ws.Range("C215").Value = 94073
$ws.Range("D215").Value = 137637989
$ws.Range("C218").Value = 50635
$ws.Range("D218").Value = 73181166
$ws.Range("C219").Value = 31
$ws.Range("D219").Value = 44422
$ws.Range("C221").Value = 4607
$ws.Range("D221").Value = 6467232
$ws.Range("C224").Value = 5545
$ws.Range("D224").Value = 7669273
$ws.Range("C227").Value = 104548
$ws.Range("D227").Value = 130877892
$ws.Range("C232").Value = 562
$ws.Range("D232").Value = 820939
$ws.Range("C234").Value = 48983
$ws.Range("D234").Value = 71766007
$ws.Range("C236").Value = 12184
$ws.Range("D236").Value = 17516277
$ws.Range("C240").Value = 2423
$ws.Range("D240").Value = 3385150
$ws.Range("C241").Value = 252961
$ws.Range("D241").Value = 319495817
$ws.Range("C243").Value = 246
$ws.Range("D243").Value = 353457
$ws.Range("C249").Value = 94651
$ws.Range("D249").Value = 138705244
$ws.Range("C252").Value = 63858
$ws.Range("D252").Value = 92542620
$ws.Range("C254").Value = 2376
$ws.Range("D254").Value = 3353361
$ws.Range("C257").Value = 4472
$ws.Range("D257").Value = 6276892
